$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12320
$ws1.Range("F3").Value = 257
$ws1.Range("F6").Value = 231
$ws1.Range("F7").Value = 12259
$ws1.Range("F8").Value = 517
$ws1.Range("F11").Value = 622
$ws1.Range("F12").Value = 2815
$ws1.Range("F13").Value = 6005
$ws1.Range("F15").Value = 3576

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12320
$ws4.Range("F3").Value = 257
$ws4.Range("F7").Value = 231
$ws4.Range("F8").Value = 12259
$ws4.Range("F9").Value = 517
$ws4.Range("F12").Value = 622
$ws4.Range("F13").Value = 2815
$ws4.Range("F15").Value = 6005
$ws4.Range("F17").Value = 3576
